{"js": "// Replace the three-digit-by-one-digit multiplication expressions with the\n// newly generated set of problems. Each original expression occurs exactly\n// once in the document body, so a scoped search-and-replace for each pair\n// is safe and unambiguous.\nconst replacements = [\n  [\"617\u00d74=2468\", \"828\u00d79=7452\"],\n  [\"657\u00d73=1971\", \"602\u00d72=1204\"],\n  [\"434\u00d78=3472\", \"261\u00d73=783\"],\n  [\"399\u00d77=2793\", \"426\u00d73=1278\"],\n  [\"573\u00d73=1719\", \"613\u00d72=1226\"],\n  [\"487\u00d72=974\", \"345\u00d74=1380\"],\n  [\"218\u00d75=1090\", \"820\u00d76=4920\"],\n  [\"429\u00d76=2574\", \"981\u00d72=1962\"],\n  [\"585\u00d72=1170\", \"316\u00d72=632\"],\n  [\"132\u00d79=1188\", \"459\u00d75=2295\"],\n  [\"448\u00d78=3584\", \"163\u00d72=326\"],\n  [\"369\u00d72=738\", \"516\u00d79=4644\"],\n  [\"972\u00d72=1944\", \"607\u00d74=2428\"],\n  [\"435\u00d74=1740\", \"904\u00d75=4520\"],\n  [\"764\u00d78=6112\", \"549\u00d79=4941\"],\n  [\"574\u00d79=5166\", \"843\u00d76=5058\"],\n  [\"276\u00d76=1656\", \"584\u00d76=3504\"],\n  [\"982\u00d73=2946\", \"419\u00d72=838\"],\n  [\"312\u00d72=624\", \"472\u00d73=1416\"],\n  [\"604\u00d77=4228\", \"540\u00d72=1080\"],\n  [\"293\u00d73=879\", \"265\u00d77=1855\"],\n  [\"965\u00d72=1930\", \"852\u00d78=6816\"],\n  [\"185\u00d72=370\", \"743\u00d79=6687\"],\n  [\"692\u00d77=4844\", \"390\u00d73=1170\"],\n  [\"311\u00d76=1866\", \"807\u00d75=4035\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit-by-one-digit multiplication expressions with the\n# newly generated set of problems. Each original expression occurs exactly\n# once in the document, so a simple Find/Replace per pair is safe.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"617\u00d74=2468\", \"828\u00d79=7452\"),\n    @(\"657\u00d73=1971\", \"602\u00d72=1204\"),\n    @(\"434\u00d78=3472\", \"261\u00d73=783\"),\n    @(\"399\u00d77=2793\", \"426\u00d73=1278\"),\n    @(\"573\u00d73=1719\", \"613\u00d72=1226\"),\n    @(\"487\u00d72=974\",  \"345\u00d74=1380\"),\n    @(\"218\u00d75=1090\", \"820\u00d76=4920\"),\n    @(\"429\u00d76=2574\", \"981\u00d72=1962\"),\n    @(\"585\u00d72=1170\", \"316\u00d72=632\"),\n    @(\"132\u00d79=1188\", \"459\u00d75=2295\"),\n    @(\"448\u00d78=3584\", \"163\u00d72=326\"),\n    @(\"369\u00d72=738\",  \"516\u00d79=4644\"),\n    @(\"972\u00d72=1944\", \"607\u00d74=2428\"),\n    @(\"435\u00d74=1740\", \"904\u00d75=4520\"),\n    @(\"764\u00d78=6112\", \"549\u00d79=4941\"),\n    @(\"574\u00d79=5166\", \"843\u00d76=5058\"),\n    @(\"276\u00d76=1656\", \"584\u00d76=3504\"),\n    @(\"982\u00d73=2946\", \"419\u00d72=838\"),\n    @(\"312\u00d72=624\",  \"472\u00d73=1416\"),\n    @(\"604\u00d77=4228\", \"540\u00d72=1080\"),\n    @(\"293\u00d73=879\",  \"265\u00d77=1855\"),\n    @(\"965\u00d72=1930\", \"852\u00d78=6816\"),\n    @(\"185\u00d72=370\",  \"743\u00d79=6687\"),\n    @(\"692\u00d77=4844\", \"390\u00d73=1170\"),\n    @(\"311\u00d76=1866\", \"807\u00d75=4035\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
